# update show by login or logout
# Adds a "type" column value of "general" to the rows on Sheet1 that did
# not yet carry a type (column M), matching the existing "private" value
# already present on the "passwd" row (M5). The new cells reuse the same
# wrap-text style already applied to the rest of column M/J/K/L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(2, 3, 4, 6, 7)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 13)   # column M
    $cell.Value = "general"
    $cell.WrapText = $true
}

# Match the author's final selection (bottom-right frozen pane) on M6.
$ws.Activate() | Out-Null
$ws.Range("M6").Select() | Out-Null
